# Auto-generated edit script applying odds updates per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("O7").Value = 1.62
$ws.Range("P7").Value = 2.2
$ws.Range("S7").Value = 6.5
$ws.Range("T7").Value = 1.11
$ws.Range("AR7").Value = 5.06
$ws.Range("AS7").Value = 1.15

# Row 8
$ws.Range("G8").Value = 2.25
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 2.88
$ws.Range("AK8").Value = 17
$ws.Range("AL8").Value = 11

# Row 9
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.95

# Row 10
$ws.Range("G10").Value = 2.75
$ws.Range("I10").Value = 2.35
$ws.Range("AF10").Value = 6.5
$ws.Range("AJ10").Value = 10

# Row 11
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 4
$ws.Range("T11").Value = 1.22
$ws.Range("AP11").Value = 1.63
$ws.Range("AQ11").Value = 2.25
$ws.Range("AR11").Value = 3.2
$ws.Range("AS11").Value = 1.34

# Row 12
$ws.Range("H12").Value = 2.87
$ws.Range("L12").Value = 4.2
$ws.Range("P12").Value = 2.4
$ws.Range("AA12").Value = 9.25
$ws.Range("AO12").Value = 55

# Row 13
$ws.Range("G13").Value = 2.65
$ws.Range("H13").Value = 2.7
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 3.3
$ws.Range("K13").Value = 1.85
$ws.Range("L13").Value = 3.7
$ws.Range("N13").Value = 5.1
$ws.Range("P13").Value = 2.35
$ws.Range("Q13").Value = 2.55
$ws.Range("W13").Value = 2.02
$ws.Range("Y13").Value = 6.4
$ws.Range("Z13").Value = 12
$ws.Range("AA13").Value = 10.25
$ws.Range("AB13").Value = 32
$ws.Range("AC13").Value = 27
$ws.Range("AE13").Value = 5.1
$ws.Range("AF13").Value = 5.4
$ws.Range("AG13").Value = 16.5
$ws.Range("AH13").Value = 100
$ws.Range("AL13").Value = 11.25
$ws.Range("AN13").Value = 32

# Row 15
$ws.Range("N15").Value = 13
$ws.Range("Q15").Value = 1.75
$ws.Range("R15").Value = 2.05

# Row 16
$ws.Range("G16").Value = 2.38
$ws.Range("I16").Value = 2.8
$ws.Range("L16").Value = 3.5
$ws.Range("Q16").Value = 1.98
$ws.Range("R16").Value = 1.88
$ws.Range("S16").Value = 3.4
$ws.Range("T16").Value = 1.3
$ws.Range("AG16").Value = 15
$ws.Range("AH16").Value = 51
$ws.Range("AK16").Value = 15
$ws.Range("AN16").Value = 23

# Row 19
$ws.Range("G19").Value = 2.35
$ws.Range("L19").Value = 3.6
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("Q19").Value = 2.08
$ws.Range("R19").Value = 1.73
$ws.Range("AB19").Value = 23
$ws.Range("AC19").Value = 21
$ws.Range("AF19").Value = 6

# Row 21
$ws.Range("G21").Value = 2.4
$ws.Range("H21").Value = 3.3
$ws.Range("I21").Value = 2.8
$ws.Range("J21").Value = 3.1
$ws.Range("K21").Value = 2.1
$ws.Range("M21").Value = 1.06
$ws.Range("N21").Value = 10
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 1.85
$ws.Range("S21").Value = 3.4
$ws.Range("T21").Value = 1.3
$ws.Range("W21").Value = 1.75
$ws.Range("X21").Value = 2
$ws.Range("AB21").Value = 23
$ws.Range("AD21").Value = 29
$ws.Range("AE21").Value = 10
$ws.Range("AH21").Value = 51
$ws.Range("AJ21").Value = 9

# Row 22
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 10
$ws.Range("O22").Value = 1.36
$ws.Range("P22").Value = 3
$ws.Range("Q22").Value = 2.1
$ws.Range("R22").Value = 1.7
$ws.Range("S22").Value = 3.75
$ws.Range("T22").Value = 1.25
$ws.Range("U22").Value = 1.44
$ws.Range("V22").Value = 2.63
$ws.Range("W22").Value = 2
$ws.Range("X22").Value = 1.75
$ws.Range("AE22").Value = 9.5
$ws.Range("AG22").Value = 19
$ws.Range("AH22").Value = 67
$ws.Range("AI22").Value = 401

# Row 24
$ws.Range("G24").Value = 1.9
$ws.Range("I24").Value = 4
$ws.Range("AP24").Value = 1.85
$ws.Range("AQ24").Value = 2

# Row 26
$ws.Range("G26").Value = 3.1
$ws.Range("H26").Value = 3.1
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 3.2
$ws.Range("M26").Value = 1.08
$ws.Range("N26").Value = 8
$ws.Range("O26").Value = 1.4
$ws.Range("P26").Value = 2.75
$ws.Range("Q26").Value = 2.25
$ws.Range("R26").Value = 1.62
$ws.Range("S26").Value = 4
$ws.Range("T26").Value = 1.22
$ws.Range("U26").Value = 1.5
$ws.Range("V26").Value = 2.5
$ws.Range("W26").Value = 1.95
$ws.Range("X26").Value = 1.8
$ws.Range("AA26").Value = 12
$ws.Range("AD26").Value = 41
$ws.Range("AE26").Value = 8
$ws.Range("AI26").Value = 351
$ws.Range("AJ26").Value = 7
$ws.Range("AL26").Value = 10
$ws.Range("AO26").Value = 34

# Row 27
$ws.Range("H27").Value = 3.4
$ws.Range("I27").Value = 3.25
$ws.Range("M27").Value = 1.04
$ws.Range("N27").Value = 12
$ws.Range("Q27").Value = 1.8
$ws.Range("R27").Value = 2
$ws.Range("U27").Value = 1.36
$ws.Range("V27").Value = 3
$ws.Range("W27").Value = 1.67
$ws.Range("X27").Value = 2.1
$ws.Range("Y27").Value = 8.5
$ws.Range("AD27").Value = 23
$ws.Range("AE27").Value = 12
$ws.Range("AI27").Value = 151
$ws.Range("AJ27").Value = 12
$ws.Range("AO27").Value = 29

# Row 28
$ws.Range("I28").Value = 3.05
$ws.Range("J28").Value = 2.75
$ws.Range("L28").Value = 3.6
$ws.Range("P28").Value = 3.4
$ws.Range("R28").Value = 1.93
$ws.Range("S28").Value = 2.57
$ws.Range("T28").Value = 1.39
$ws.Range("Y28").Value = 9.75
$ws.Range("Z28").Value = 13
$ws.Range("AB28").Value = 24
$ws.Range("AC28").Value = 16.5
$ws.Range("AD28").Value = 22
$ws.Range("AE28").Value = 11
$ws.Range("AG28").Value = 11.5
$ws.Range("AJ28").Value = 10.5
$ws.Range("AK28").Value = 17
$ws.Range("AN28").Value = 25
$ws.Range("AO28").Value = 29

# Row 29
$ws.Range("G29").Value = 1.45
$ws.Range("I29").Value = 7
$ws.Range("J29").Value = 1.95
$ws.Range("M29").Value = 1.03
$ws.Range("O29").Value = 1.19
$ws.Range("T29").Value = 1.37
$ws.Range("W29").Value = 1.95
$ws.Range("X29").Value = 1.8
$ws.Range("Y29").Value = 7
$ws.Range("AE29").Value = 12
$ws.Range("AG29").Value = 19
$ws.Range("AL29").Value = 21
$ws.Range("AM29").Value = 81
$ws.Range("AO29").Value = 51

# Row 30
$ws.Range("H30").Value = 3.25
$ws.Range("M30").Value = 1.04
$ws.Range("O30").Value = 1.27
$ws.Range("R30").Value = 1.75
$ws.Range("T30").Value = 1.25
$ws.Range("W30").Value = 1.8
$ws.Range("X30").Value = 1.95
$ws.Range("Y30").Value = 8

# Row 31
$ws.Range("M31").Value = 1.07
$ws.Range("O31").Value = 1.41
$ws.Range("P31").Value = 2.62
$ws.Range("T31").Value = 1.15
$ws.Range("AI31").Value = 1250
$ws.Range("AP31").Value = 1.8
$ws.Range("AQ31").Value = 2

# Row 32
$ws.Range("G32").Value = 2.8
$ws.Range("I32").Value = 2.63
$ws.Range("J32").Value = 3.6
$ws.Range("K32").Value = 1.91
$ws.Range("L32").Value = 3.5
$ws.Range("M32").Value = 1.11
$ws.Range("N32").Value = 6.5
$ws.Range("W32").Value = 2.1
$ws.Range("X32").Value = 1.67
$ws.Range("Y32").Value = 7
$ws.Range("Z32").Value = 13
$ws.Range("AB32").Value = 29
$ws.Range("AC32").Value = 26
$ws.Range("AJ32").Value = 7
$ws.Range("AK32").Value = 12
$ws.Range("AL32").Value = 11
$ws.Range("AM32").Value = 26
$ws.Range("AN32").Value = 26
